$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6671080
$ws.Range("H41").Value = 805
$ws.Range("J41").Value = 958.2
$ws.Range("L41").Value = 958.2
$ws.Range("N41").Value = -1838.2
$ws.Range("H43").Value = 7536.769
$ws.Range("I43").Value = 4000
$ws.Range("J43").Value = 9747.25
$ws.Range("K43").Value = 4000
$ws.Range("L43").Value = 9747.25
$ws.Range("M43").Value = -3931
$ws.Range("N43").Value = -9885.25
$ws.Range("H86").Value = 3888.3
$ws.Range("I86").Value = 2579.6
$ws.Range("K86").Value = 2579.6
$ws.Range("M86").Value = -1456.6
$ws.Range("H89").Value = 3888.3
$ws.Range("I89").Value = 2579.6
$ws.Range("K89").Value = 12898
$ws.Range("M89").Value = -7282
$ws.Range("H106").Value = 107334.3
$ws.Range("I106").Value = 107334.3
$ws.Range("K106").Value = 107334.3
$ws.Range("M106").Value = -106703.3
$ws.Range("H111").Value = 41976.875
$ws.Range("I111").Value = 12433.077
$ws.Range("J111").Value = 170000
$ws.Range("K111").Value = 37299.231
$ws.Range("L111").Value = 510000
$ws.Range("M111").Value = -34232.231
$ws.Range("N111").Value = -516134
$ws.Range("H132").Value = 31253582
$ws.Range("I132").Value = 33337086
$ws.Range("J132").Value = 1003
$ws.Range("K132").Value = 100011258
$ws.Range("L132").Value = 3009
$ws.Range("M132").Value = -100008728
$ws.Range("N132").Value = -8069
$ws.Range("H135").Value = 4022.8262
$ws.Range("I135").Value = 1197.4286
$ws.Range("J135").Value = 8417.888999999999
$ws.Range("K135").Value = 10776.8574
$ws.Range("L135").Value = 75761.00099999999
$ws.Range("M135").Value = -8241.857399999999
$ws.Range("N135").Value = -80831.00099999999
$ws.Range("H138").Value = 1583657.2
$ws.Range("J138").Value = 2002461.5
$ws.Range("L138").Value = 6007384.5
$ws.Range("N138").Value = -6017664.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 1247.5
$ws.Range("J17").Value = 1990
$ws.Range("L17").Value = 1990
$ws.Range("N17").Value = -2336
$ws.Range("H32").Value = 1244.7073
$ws.Range("I32").Value = 889.1818
$ws.Range("K32").Value = 889.1818
$ws.Range("M32").Value = -602.1818
$ws.Range("H45").Value = 2244.6924
$ws.Range("I45").Value = 1449.7142
$ws.Range("K45").Value = 1449.7142
$ws.Range("M45").Value = -1072.7142
$ws.Range("H61").Value = 4688.2925
$ws.Range("I61").Value = 2792.2258
$ws.Range("K61").Value = 2792.2258
$ws.Range("M61").Value = -2580.2258
$ws.Range("H110").Value = 2323.182
$ws.Range("I110").Value = 1615.8889
$ws.Range("K110").Value = 1615.8889
$ws.Range("M110").Value = 429.1111000000001
$ws.Range("H132").Value = 3276.2144
$ws.Range("I132").Value = 2469.2092
$ws.Range("K132").Value = 7407.6276
$ws.Range("M132").Value = -4877.6276
$ws.Range("H136").Value = 4688.2925
$ws.Range("I136").Value = 2792.2258
$ws.Range("K136").Value = 8376.6774
$ws.Range("M136").Value = -5826.6774

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3477.9048
$ws.Range("I105").Value = 1945.5625
$ws.Range("K105").Value = 1945.5625
$ws.Range("M105").Value = -198.5625
$ws.Range("H107").Value = 2325.75
$ws.Range("I107").Value = 2472.3333
$ws.Range("K107").Value = 2472.3333
$ws.Range("M107").Value = -552.3332999999998

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2405.7026
$ws.Range("I31").Value = 1028.64
$ws.Range("J31").Value = 5274.5835
$ws.Range("K31").Value = 1028.64
$ws.Range("L31").Value = 5274.5835
$ws.Range("M31").Value = -733.6400000000001
$ws.Range("N31").Value = -5864.5835
$ws.Range("H34").Value = 2405.7026
$ws.Range("I34").Value = 1028.64
$ws.Range("J34").Value = 5274.5835
$ws.Range("K34").Value = 1028.64
$ws.Range("L34").Value = 5274.5835
$ws.Range("M34").Value = -826.6400000000001
$ws.Range("N34").Value = -5678.5835
$ws.Range("H105").Value = 3280.75
$ws.Range("I105").Value = 2624.5
$ws.Range("K105").Value = 2624.5
$ws.Range("M105").Value = -877.5
$ws.Range("H107").Value = 523.05884
$ws.Range("J107").Value = 534.125
$ws.Range("L107").Value = 534.125
$ws.Range("N107").Value = -4374.125
$ws.Range("H134").Value = 1637.5385
$ws.Range("I134").Value = 1595.8182
$ws.Range("K134").Value = 4787.4546
$ws.Range("M134").Value = -2252.4546

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2019.6
$ws.Range("I68").Value = 1933.3334
$ws.Range("J68").Value = 2149
$ws.Range("K68").Value = 5800.0002
$ws.Range("L68").Value = 6447
$ws.Range("M68").Value = -4989.0002
$ws.Range("N68").Value = -8069
$ws.Range("H71").Value = 2019.6
$ws.Range("I71").Value = 1933.3334
$ws.Range("J71").Value = 2149
$ws.Range("K71").Value = 17400.0006
$ws.Range("L71").Value = 19341
$ws.Range("M71").Value = -13344.0006
$ws.Range("N71").Value = -27453
$ws.Range("H92").Value = 1125
$ws.Range("J92").Value = 1125
$ws.Range("L92").Value = 3375
$ws.Range("N92").Value = -5871

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3836.7778
$ws.Range("I80").Value = 3205
$ws.Range("K80").Value = 3205
$ws.Range("M80").Value = -2207
$ws.Range("H83").Value = 3836.7778
$ws.Range("I83").Value = 3205
$ws.Range("K83").Value = 16025
$ws.Range("M83").Value = -11033
$ws.Range("J122").Value = 4499
$ws.Range("L122").Value = 13497
$ws.Range("N122").Value = -18397
$ws.Range("H132").Value = 3495.3635
$ws.Range("I132").Value = 3876.9092
$ws.Range("J132").Value = 1969.1818
$ws.Range("K132").Value = 11630.7276
$ws.Range("L132").Value = 5907.5454
$ws.Range("M132").Value = -9100.7276
$ws.Range("N132").Value = -10967.5454

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1498.1904
$ws.Range("I16").Value = 1583.8889
$ws.Range("K16").Value = 1583.8889
$ws.Range("M16").Value = -1413.8889
$ws.Range("H61").Value = 14631.148
$ws.Range("I61").Value = 10978.579
$ws.Range("J61").Value = 23306
$ws.Range("K61").Value = 10978.579
$ws.Range("L61").Value = 23306
$ws.Range("M61").Value = -10776.579
$ws.Range("N61").Value = -23710
$ws.Range("H82").Value = 8457.75
$ws.Range("I82").Value = 10393.583
$ws.Range("K82").Value = 10393.583
$ws.Range("M82").Value = -10032.583
$ws.Range("H85").Value = 8457.75
$ws.Range("I85").Value = 10393.583
$ws.Range("K85").Value = 10393.583
$ws.Range("M85").Value = -9145.583000000001
$ws.Range("H113").Value = 14631.148
$ws.Range("I113").Value = 10978.579
$ws.Range("J113").Value = 23306
$ws.Range("K113").Value = 10978.579
$ws.Range("L113").Value = 23306
$ws.Range("M113").Value = -8808.579
$ws.Range("N113").Value = -27646
$ws.Range("H132").Value = 3799
$ws.Range("I132").Value = 3799
$ws.Range("K132").Value = 11397
$ws.Range("M132").Value = -8867

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 160332.92
$ws.Range("I2").Value = 7416.125
$ws.Range("K2").Value = 7416.125
$ws.Range("M2").Value = -7304.125
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("H74").Value = 17089.092
$ws.Range("J74").Value = 15140.571
$ws.Range("L74").Value = 15140.571
$ws.Range("N74").Value = -17012.571
$ws.Range("H77").Value = 17089.092
$ws.Range("J77").Value = 15140.571
$ws.Range("L77").Value = 45421.713
$ws.Range("N77").Value = -54781.713
$ws.Range("H113").Value = 833.13336
$ws.Range("I113").Value = 608.5625
$ws.Range("J113").Value = 1089.7858
$ws.Range("K113").Value = 1825.6875
$ws.Range("L113").Value = 3269.3574
$ws.Range("M113").Value = 344.3125
$ws.Range("N113").Value = -7609.357400000001
$ws.Range("H132").Value = 2250
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("M132").Value = -3470
$ws.Range("H136").Value = 5733.305
$ws.Range("I136").Value = 4976.548
$ws.Range("K136").Value = 14929.644
$ws.Range("M136").Value = -12379.644

# Special case: clear M40 on WVR (cell removed from source entirely)
$wsWVR = $wb.Worksheets.Item("WVR")
$wsWVR.Range("M40").ClearContents()
